{"js": "// Replace the date line and each \"a\u00f7b=c, d\" answer cell in document order.\n// Each (old, new) pair is unique at the time it is applied (verified offline),\n// so performing them sequentially in this order never matches more than one\n// range in the document.\nconst replacements = [\n  [\"2024-08-05 Monday\", \"2024-08-06 Tuesday\"],\n  [\"19\u00f77=2, 5\", \"31\u00f77=4, 3\"],\n  [\"55\u00f79=6, 1\", \"75\u00f73=25, 0\"],\n  [\"27\u00f74=6, 3\", \"79\u00f75=15, 4\"],\n  [\"61\u00f72=30, 1\", \"49\u00f78=6, 1\"],\n  [\"31\u00f76=5, 1\", \"37\u00f78=4, 5\"],\n  [\"67\u00f79=7, 4\", \"65\u00f72=32, 1\"],\n  [\"77\u00f72=38, 1\", \"89\u00f78=11, 1\"],\n  [\"99\u00f78=12, 3\", \"88\u00f75=17, 3\"],\n  [\"66\u00f79=7, 3\", \"24\u00f77=3, 3\"],\n  [\"34\u00f73=11, 1\", \"87\u00f79=9, 6\"],\n  [\"64\u00f79=7, 1\", \"96\u00f73=32, 0\"],\n  [\"81\u00f77=11, 4\", \"72\u00f75=14, 2\"],\n  [\"40\u00f76=6, 4\", \"87\u00f78=10, 7\"],\n  [\"28\u00f75=5, 3\", \"77\u00f74=19, 1\"],\n  [\"57\u00f76=9, 3\", \"81\u00f79=9, 0\"],\n  [\"26\u00f78=3, 2\", \"15\u00f79=1, 6\"],\n  [\"88\u00f79=9, 7\", \"56\u00f79=6, 2\"],\n  [\"80\u00f76=13, 2\", \"17\u00f75=3, 2\"],\n  [\"46\u00f77=6, 4\", \"22\u00f74=5, 2\"],\n  [\"96\u00f78=12, 0\", \"89\u00f73=29, 2\"],\n  [\"84\u00f73=28, 0\", \"76\u00f73=25, 1\"],\n  [\"11\u00f74=2, 3\", \"77\u00f72=38, 1\"],\n  [\"14\u00f77=2, 0\", \"33\u00f74=8, 1\"],\n  [\"44\u00f78=5, 4\", \"58\u00f79=6, 4\"],\n  [\"90\u00f77=12, 6\", \"68\u00f74=17, 0\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replace only the first (and expected only) match to stay safe even if\n  // duplicates were ever introduced.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Sequentially replace the date line and each answer cell's text in\n# document order, using Word's Find/Replace (wdReplaceAll = 2, wdFindContinue = 1).\n# Each (old, new) pair is unique in the live document at the moment it runs\n# (verified offline), so every call here matches exactly one occurrence.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-05 Monday\", \"2024-08-06 Tuesday\"),\n    @(\"19\u00f77=2, 5\", \"31\u00f77=4, 3\"),\n    @(\"55\u00f79=6, 1\", \"75\u00f73=25, 0\"),\n    @(\"27\u00f74=6, 3\", \"79\u00f75=15, 4\"),\n    @(\"61\u00f72=30, 1\", \"49\u00f78=6, 1\"),\n    @(\"31\u00f76=5, 1\", \"37\u00f78=4, 5\"),\n    @(\"67\u00f79=7, 4\", \"65\u00f72=32, 1\"),\n    @(\"77\u00f72=38, 1\", \"89\u00f78=11, 1\"),\n    @(\"99\u00f78=12, 3\", \"88\u00f75=17, 3\"),\n    @(\"66\u00f79=7, 3\", \"24\u00f77=3, 3\"),\n    @(\"34\u00f73=11, 1\", \"87\u00f79=9, 6\"),\n    @(\"64\u00f79=7, 1\", \"96\u00f73=32, 0\"),\n    @(\"81\u00f77=11, 4\", \"72\u00f75=14, 2\"),\n    @(\"40\u00f76=6, 4\", \"87\u00f78=10, 7\"),\n    @(\"28\u00f75=5, 3\", \"77\u00f74=19, 1\"),\n    @(\"57\u00f76=9, 3\", \"81\u00f79=9, 0\"),\n    @(\"26\u00f78=3, 2\", \"15\u00f79=1, 6\"),\n    @(\"88\u00f79=9, 7\", \"56\u00f79=6, 2\"),\n    @(\"80\u00f76=13, 2\", \"17\u00f75=3, 2\"),\n    @(\"46\u00f77=6, 4\", \"22\u00f74=5, 2\"),\n    @(\"96\u00f78=12, 0\", \"89\u00f73=29, 2\"),\n    @(\"84\u00f73=28, 0\", \"76\u00f73=25, 1\"),\n    @(\"11\u00f74=2, 3\", \"77\u00f72=38, 1\"),\n    @(\"14\u00f77=2, 0\", \"33\u00f74=8, 1\"),\n    @(\"44\u00f78=5, 4\", \"58\u00f79=6, 4\"),\n    @(\"90\u00f77=12, 6\", \"68\u00f74=17, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\n"}
